$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("S6").Value = 1664

$textRows = @{
    8  = "146"
    9  = "78"
    10 = "149"
    11 = "82"
    12 = "152"
    13 = "86"
    14 = "155"
    15 = "90"
    16 = "158"
    17 = "94"
    18 = "161"
    19 = "121"
    20 = "164"
    21 = "167"
    22 = "170"
    24 = "173"
    25 = "140"
    26 = "176"
}

foreach ($row in $textRows.Keys) {
    $cell = $ws.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $textRows[$row]
    $cell.Style = "Normal"
}

$ws.Range("S32").Value = 28307

$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
